{"js": "const replacements = [\n  [\"2024-07-05 Friday\", \"2024-07-06 Saturday\"],\n  [\"619\u00f73=\", \"583\u00f75=\"],\n  [\"943\u00f76=\", \"298\u00f79=\"],\n  [\"925\u00f78=\", \"487\u00f77=\"],\n  [\"976\u00f72=\", \"292\u00f75=\"],\n  [\"612\u00f73=\", \"207\u00f77=\"],\n  [\"881\u00f79=\", \"364\u00f79=\"],\n  [\"684\u00f72=\", \"235\u00f79=\"],\n  [\"376\u00f75=\", \"528\u00f74=\"],\n  [\"163\u00f79=\", \"397\u00f77=\"],\n  [\"334\u00f72=\", \"604\u00f77=\"],\n  [\"713\u00f75=\", \"649\u00f78=\"],\n  [\"365\u00f78=\", \"575\u00f79=\"],\n  [\"408\u00f73=\", \"368\u00f76=\"],\n  [\"971\u00f73=\", \"176\u00f77=\"],\n  [\"171\u00f74=\", \"915\u00f79=\"],\n  [\"725\u00f77=\", \"426\u00f78=\"],\n  [\"664\u00f77=\", \"520\u00f78=\"],\n  [\"655\u00f79=\", \"386\u00f75=\"],\n  [\"858\u00f73=\", \"908\u00f78=\"],\n  [\"400\u00f72=\", \"779\u00f75=\"],\n  [\"658\u00f78=\", \"653\u00f74=\"],\n  [\"153\u00f72=\", \"880\u00f78=\"],\n  [\"631\u00f76=\", \"505\u00f77=\"],\n  [\"268\u00f76=\", \"901\u00f72=\"],\n  [\"940\u00f78=\", \"631\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-05 Friday\", \"2024-07-06 Saturday\"),\n    @(\"619\u00f73=\", \"583\u00f75=\"),\n    @(\"943\u00f76=\", \"298\u00f79=\"),\n    @(\"925\u00f78=\", \"487\u00f77=\"),\n    @(\"976\u00f72=\", \"292\u00f75=\"),\n    @(\"612\u00f73=\", \"207\u00f77=\"),\n    @(\"881\u00f79=\", \"364\u00f79=\"),\n    @(\"684\u00f72=\", \"235\u00f79=\"),\n    @(\"376\u00f75=\", \"528\u00f74=\"),\n    @(\"163\u00f79=\", \"397\u00f77=\"),\n    @(\"334\u00f72=\", \"604\u00f77=\"),\n    @(\"713\u00f75=\", \"649\u00f78=\"),\n    @(\"365\u00f78=\", \"575\u00f79=\"),\n    @(\"408\u00f73=\", \"368\u00f76=\"),\n    @(\"971\u00f73=\", \"176\u00f77=\"),\n    @(\"171\u00f74=\", \"915\u00f79=\"),\n    @(\"725\u00f77=\", \"426\u00f78=\"),\n    @(\"664\u00f77=\", \"520\u00f78=\"),\n    @(\"655\u00f79=\", \"386\u00f75=\"),\n    @(\"858\u00f73=\", \"908\u00f78=\"),\n    @(\"400\u00f72=\", \"779\u00f75=\"),\n    @(\"658\u00f78=\", \"653\u00f74=\"),\n    @(\"153\u00f72=\", \"880\u00f78=\"),\n    @(\"631\u00f76=\", \"505\u00f77=\"),\n    @(\"268\u00f76=\", \"901\u00f72=\"),\n    @(\"940\u00f78=\", \"631\u00f73=\"),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}"}
